# feat: add 2022-Q3 data
#
# The workbook tracks quarterly fund-holding snapshots, one sheet per
# quarter, plus a "总计" (totals) roll-up sheet. This change adds a new
# "2022-Q3" quarter sheet (inserted right after the roll-up sheet and
# before the existing "2022-Q2" sheet) and records its summary row in
# "总计". The pre-existing quarter sheets ("2022-Q2", "2022-Q1",
# "2021-Q4") are left untouched.

$wb = $excel.ActiveWorkbook

# Helper: force a cell to be written as TEXT (so things like fund codes
# "004854" or percentages "15.34" keep their original string shape
# instead of being coerced to numbers), while keeping the cell's visual
# style identical to a neighboring "plain" cell (formatDonor) instead of
# the blended "text number format" style that NumberFormat="@" would
# otherwise leave behind.
function Set-TextValue($range, $value, $formatDonor) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $formatDonor.Copy()
    $range.PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q3" sheet, positioned before "2022-Q2".
#    Copying "2022-Q2" gives us the exact same header row / column
#    styles for free; we then overwrite the data with the new quarter's
#    numbers.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Row 2: 广发中证全指汽车指数A (004854)
$q3.Range("A2").Value = 0
Set-TextValue $q3.Range("B2") "004854" $q3.Range("C2")
$q3.Range("C2").Value = "广发中证全指汽车指数A"
Set-TextValue $q3.Range("D2") "15.34" $q3.Range("C2")
Set-TextValue $q3.Range("E2") "94.07" $q3.Range("C2")
Set-TextValue $q3.Range("F2") "2.45" $q3.Range("C2")
Set-TextValue $q3.Range("G2") "0.3758" $q3.Range("C2")
$q3.Range("H2").Value = 9

# Row 3: 广发中证全指汽车指数C (004855)
$q3.Range("A3").Value = 1
Set-TextValue $q3.Range("B3") "004855" $q3.Range("C3")
$q3.Range("C3").Value = "广发中证全指汽车指数C"
Set-TextValue $q3.Range("D3") "10.15" $q3.Range("C3")
Set-TextValue $q3.Range("E3") "94.07" $q3.Range("C3")
Set-TextValue $q3.Range("F3") "2.45" $q3.Range("C3")
Set-TextValue $q3.Range("G3") "0.2487" $q3.Range("C3")
$q3.Range("H3").Value = 9

# ---------------------------------------------------------------------
# 2) Update the "总计" roll-up sheet: shift the existing quarters down
#    one row, relabel them, and append the new 2021-Q4 row that falls
#    out the bottom, plus fill in the new 2022-Q3 total at the top.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.62

$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 1

$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.96

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q4"
$total.Range("C5").Value = 2
$total.Range("D5").Value = 1.04

# New row 5 needs the same style as the rows above it (bordered, bold,
# centered) since it didn't exist before.
$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)
$total.Range("A5").Value = 3
